# Regenerate the "K" column (column G) values for rows 2-20 on Sheet1.
# This mirrors the source data regeneration described in the commit message
# ("regen save_data to use K instead of Strike#, regen std/mean, calc and
# write s_vals") — here it manifests purely as updated literal values in
# column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 6
    3  = 2
    4  = 3
    5  = 10
    6  = 8
    7  = 5
    8  = 5
    9  = 6
    10 = 6
    11 = 5
    12 = 4
    13 = 6
    14 = 5
    15 = 2
    16 = 7
    17 = 7
    18 = 8
    19 = 4
    20 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
